{"js": "// \"Burrower\" trait was renamed to \"Burrowing\" and its rules text was\n// rewritten to give distance-based burrowing speeds instead of a flat\n// \"submerge\" description.\n//\n// Old:  Burrower. The radscorpion has a burrowing speed of 5 feet, which\n//       it can use to fully submerge itself underneath loose rock, sand,\n//       or dirt.\n// New:  Burrowing. The radscorpion has a burrowing speed of 15 feet\n//       through loose earth, 10 feet through solid rock, and 0 feet\n//       through metal.\n\n// 1) Rename the bold trait label, leaving its bold formatting untouched.\nconst label = context.document.body.search(\"Burrower.\", { matchCase: true });\nlabel.load(\"items\");\nawait context.sync();\n\nif (label.items.length > 0) {\n  label.items[0].insertText(\"Burrowing.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Replace the rest of the sentence with the new burrowing-speed rules.\nconst oldSentence =\n  \"The radscorpion has a burrowing speed of 5 feet, which it can use to \" +\n  \"fully submerge itself underneath loose rock, sand, or dirt.\";\nconst newSentence =\n  \"The radscorpion has a burrowing speed of 15 feet through loose earth, \" +\n  \"10 feet through solid rock, and 0 feet through metal.\";\n\nconst body = context.document.body.search(oldSentence, { matchCase: true });\nbody.load(\"items\");\nawait context.sync();\n\nif (body.items.length > 0) {\n  body.items[0].insertText(newSentence, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Burrower\" trait was renamed to \"Burrowing\" and its rules text was\n# rewritten to give distance-based burrowing speeds instead of a flat\n# \"submerge\" description.\n#\n# Old:  Burrower. The radscorpion has a burrowing speed of 5 feet, which\n#       it can use to fully submerge itself underneath loose rock, sand,\n#       or dirt.\n# New:  Burrowing. The radscorpion has a burrowing speed of 15 feet\n#       through loose earth, 10 feet through solid rock, and 0 feet\n#       through metal.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n# 1) Rename the bold trait label, leaving its bold formatting untouched.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Burrower.\"\n$find.Replacement.Text = \"Burrowing.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n\n# 2) Replace the rest of the sentence with the new burrowing-speed rules.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"The radscorpion has a burrowing speed of 5 feet, which it can use to fully submerge itself underneath loose rock, sand, or dirt.\"\n$find2.Replacement.Text = \"The radscorpion has a burrowing speed of 15 feet through loose earth, 10 feet through solid rock, and 0 feet through metal.\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceAll)\n"}
